# Apply the cryptos-list price/volume refresh described by the commit.
# Column D ("Price") and E ("Volume(1h)") are plain text cells in the source
# workbook (inline strings) -- many of the new Price values look like plain
# decimals (e.g. "310.50"), which Excel would otherwise auto-convert to a
# Number (dropping the trailing zero / changing the stored type). To preserve
# the original text semantics we write those through the classic
# apostrophe-prefix "force text" entry method and then reset the cell style so
# no stray NumberFormat/quotePrefix is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.492.14"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.365.86"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'310.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'104.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.27%  "
$ws.Range("D7").Value = "'0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").Value = "'36.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "'52.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "'0.0815"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "2.732.32"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "'15.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("D17").Value = "2.360.75"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "'0.815"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "43.476.56"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "'12.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.52%  "
$ws.Range("D21").Value = "0.0₃0931"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'6.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.81%  "
$ws.Range("D23").Value = "'68.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'243.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Value = "'26.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.80%  "
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "'36.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'162.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'5.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'18.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'2.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.46%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.01%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0742"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  +10.53%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'2.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("D44").Value = "'20.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").Value = "2.007.07"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "'3.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "'10.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.46%  "
$ws.Range("D49").Value = "'58.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("E51").Value = "  +3.33%  "
